$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force these D-column cells to Text format so numeric-looking price strings
# (e.g. "219.05") are preserved as literal text rather than being parsed as numbers.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "26.115.22"
$ws.Range("E2").Value = "  -0.62%  "
$ws.Range("D3").Value = "1.654.70"
$ws.Range("E3").Value = "  -0.75%  "
$ws.Range("E4").Value = "  -0.44%  "
$ws.Range("D5").Value = "219.05"
$ws.Range("E5").Value = "  -0.35%  "
$ws.Range("D6").Value = "0.5261"
$ws.Range("E6").Value = "  -0.57%  "
$ws.Range("E7").Value = "  -0.45%  "
$ws.Range("D8").Value = "0.2683"
$ws.Range("E8").Value = "  +1.34%  "
$ws.Range("D9").Value = "0.06380"
$ws.Range("E9").Value = "  +0.20%  "
$ws.Range("D10").Value = "20.57"
$ws.Range("E10").Value = "  -1.73%  "
$ws.Range("D11").Value = "0.07688"
$ws.Range("E11").Value = "  -1.88%  "
$ws.Range("D12").Value = "4.614"
$ws.Range("E12").Value = "  +2.04%  "
$ws.Range("D13").Value = "1.721.75"
$ws.Range("E13").Value = "  +3.16%  "
$ws.Range("D14").Value = "1.883.57"
$ws.Range("E14").Value = "  -0.65%  "
$ws.Range("D15").Value = "0.5634"
$ws.Range("D16").Value = "0.0₅8250"
$ws.Range("E16").Value = "  +1.68%  "
$ws.Range("D17").Value = "65.67"
$ws.Range("E17").Value = "  -0.02%  "
$ws.Range("D18").Value = "26.109.29"
$ws.Range("E18").Value = "  -0.70%  "
$ws.Range("E19").Value = "  -0.47%  "
$ws.Range("E20").Value = "  -0.61%  "
$ws.Range("D21").Value = "10.35"
$ws.Range("E21").Value = "  +0.79%  "
$ws.Range("D22").Value = "190.27"
$ws.Range("E22").Value = "  -4.62%  "
$ws.Range("D23").Value = "5.993"
$ws.Range("E23").Value = "  -1.00%  "
$ws.Range("E24").Value = "  -0.49%  "
$ws.Range("D25").Value = "146.64"
$ws.Range("E25").Value = "  +0.44%  "
$ws.Range("E26").Value = "  -1.03%  "
$ws.Range("D27").Value = "7.261"
$ws.Range("E27").Value = "  +0.32%  "
$ws.Range("D28").Value = "16.00"
$ws.Range("E28").Value = "  -1.02%  "
$ws.Range("D29").Value = "1.523"
$ws.Range("E29").Value = "  -0.72%  "
$ws.Range("D30").Value = "0.05651"
$ws.Range("E30").Value = "  -4.37%  "
$ws.Range("D32").Value = "3.497"
$ws.Range("E32").Value = "  -0.56%  "
$ws.Range("D33").Value = "3.384"
$ws.Range("E33").Value = "  +1.77%  "
$ws.Range("D34").Value = "1.582"
$ws.Range("E34").Value = "  -0.97%  "
$ws.Range("D35").Value = "2.800"
$ws.Range("E35").Value = "  -0.77%  "
$ws.Range("E36").Value = "  -1.27%  "
$ws.Range("D37").Value = "2.410"
$ws.Range("E38").Value = "  -0.20%  "
$ws.Range("D39").Value = "0.01598"
$ws.Range("E39").Value = "  -0.92%  "
$ws.Range("D40").Value = "5.972"
$ws.Range("E40").Value = "  +0.25%  "
$ws.Range("D41").Value = "1.004"
$ws.Range("E41").Value = "  -0.49%  "
$ws.Range("D42").Value = "0.8364"
$ws.Range("E42").Value = "  -2.50%  "
$ws.Range("D43").Value = "1.022.69"
$ws.Range("E43").Value = "  -4.85%  "
$ws.Range("D44").Value = "101.31"
$ws.Range("E44").Value = "  -1.41%  "
$ws.Range("D45").Value = "1.793.69"
$ws.Range("E45").Value = "  -0.68%  "
$ws.Range("D46").Value = "58.35"
$ws.Range("E46").Value = "  -0.23%  "
$ws.Range("D47").Value = "0.0₈106"
$ws.Range("E47").Value = "  +1.04%  "
$ws.Range("E48").Value = "  -0.96%  "
$ws.Range("D49").Value = "0.05334"
$ws.Range("E49").Value = "  +3.71%  "
$ws.Range("D50").Value = "8.031"
$ws.Range("E50").Value = "  -0.29%  "
$ws.Range("D51").Value = "0.4343"
$ws.Range("E51").Value = "  -1.55%  "
